$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.875.95'
$ws.Cells.Item(2, 5).Value = '  +2.95%  '

$ws.Cells.Item(3, 4).Value = '3.034.82'
$ws.Cells.Item(3, 5).Value = '  +2.08%  '

$ws.Cells.Item(4, 5).Value = '  +0.11%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '595.44'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.19%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '152.02'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +7.04%  '

$ws.Cells.Item(7, 5).Value = '  -0.08%  '

$ws.Cells.Item(8, 4).Value = '3.029.31'
$ws.Cells.Item(8, 5).Value = '  +1.98%  '

$ws.Cells.Item(9, 5).Value = '  +0.73%  '

$ws.Cells.Item(10, 5).Value = '  +10.32%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.151'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +5.87%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '0.462'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +2.28%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '0.0000234'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +4.07%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '34.89'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +2.78%  '

$ws.Cells.Item(16, 4).Value = '3.540.11'
$ws.Cells.Item(16, 5).Value = '  +2.21%  '

$ws.Cells.Item(17, 4).Value = '62.882.66'
$ws.Cells.Item(17, 5).Value = '  +2.99%  '

$ws.Cells.Item(18, 5).Value = '  +0.70%  '

$ws.Cells.Item(19, 4).Value = '3.043.01'
$ws.Cells.Item(19, 5).Value = '  +2.53%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '456.24'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +2.01%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '14.21'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +2.49%  '

$ws.Cells.Item(22, 5).Value = '  +1.72%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '7.49'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +2.09%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '83.05'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +2.23%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +5.24%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '10.83'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +9.43%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '12.14'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.73%  '

$ws.Cells.Item(28, 5).Value = '  -0.08%  '

$ws.Cells.Item(29, 2).Value = 'PancakeSwap'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.71'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +2.62%  '

$ws.Cells.Item(30, 2).Value = 'NEARProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '7.41'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +8.78%  '

$ws.Cells.Item(31, 5).Value = '  +0.20%  '

$ws.Cells.Item(32, 5).Value = '  +5.29%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '27.60'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +1.61%  '

$ws.Cells.Item(34, 5).Value = '  +4.66%  '

$ws.Cells.Item(35, 4).Value = '0.0₃0853'
$ws.Cells.Item(35, 5).Value = '  +8.86%  '

$ws.Cells.Item(36, 5).Value = '  +2.83%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '5.90'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +3.45%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '3.11'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +13.11%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.72%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '50.46'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.51%  '

$ws.Cells.Item(41, 5).Value = '  +0.04%  '

$ws.Cells.Item(42, 5).Value = '  +4.61%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.293'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +11.83%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '41.31'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +11.16%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '392.04'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.64%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.0356'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +1.43%  '

$ws.Cells.Item(47, 4).Value = '2.744.23'
$ws.Cells.Item(47, 5).Value = '  +1.94%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '132.60'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.65%  '

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '2.20'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +2.90%  '

$ws.Cells.Item(51, 5).Value = '  +0.66%  '
